$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '319.83'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '7.51%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '49.03'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '17.37%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.274'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '5.08%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08103'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '7.68%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.606'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5.26%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.662'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.79%'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.194'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '30.05%'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1315'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '11.59%'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1942'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '6.25%'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09533'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '6.56%'

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04513'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '10.41%'

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.11%'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001330'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '3.59%'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005984'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.00%'

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.365'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.72%'

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.436'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.46%'

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.3392'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.92%'

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'MCDex'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.176'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.56%'

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1422'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.64%'

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3058'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-5.08%'

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04300'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '5.05%'

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001310'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.42%'

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004245'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '8.11%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001352'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '3.86%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003544'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-4.83%'

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02681'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '11.46%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05577'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '7.15%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006309'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.01%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007698'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.48%'

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '8.59%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007698'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.03%'

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '13.69%'

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.75%'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006997'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.19%'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.01%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06117'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '34.64%'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004005'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-4.73%'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.01%'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.01%'
